$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113; all existing rows 113..162 shift down to 114..163
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new data record
$ws.Range("A113").Value = 10
$ws.Range("B113").Value = "Vega Modelo de Temuco"
$ws.Range("C113").Value = "La Araucanía"
$ws.Range("D113").Value = 44726
$ws.Range("E113").Value = 9
$ws.Range("F113").Value = 100112013
$ws.Range("G113").Value = "Alcachofa"
$ws.Range("H113").Value = "Madrigal"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 55
$ws.Range("K113").Value = 21000
$ws.Range("L113").Value = 21000
$ws.Range("M113").Value = 21000
$ws.Range("N113").Value = "$/caja 30 unidades"
$ws.Range("O113").Value = "Provincia de Limarí"
$ws.Range("P113").Value = 700
$ws.Range("Q113").Value = 30
$ws.Range("R113").Value = "Hortaliza"
